$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update quantity (column D) values for rows 2-19
$ws.Range("D2").Value = 6
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 12
$ws.Range("D7").Value = 3
$ws.Range("D8").Value = 2
$ws.Range("D9").Value = 50
$ws.Range("D10").Value = 100
$ws.Range("D11").Value = 3
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 15
$ws.Range("D14").Value = 80
$ws.Range("D15").Value = 100
$ws.Range("D16").Value = 30
$ws.Range("D17").Value = 3
$ws.Range("D18").Value = 10
$ws.Range("D19").Value = 10

# Update shipment_id (column B) and product_id (column C) for rows 15 and 19
$ws.Range("B15").Value = 6
$ws.Range("B19").Value = 7
$ws.Range("C19").Value = 3

# Remove row 20 entirely (data previously held cargo_id 19)
$ws.Rows.Item(20).Delete()

# Update the selected cell to match the saved selection state
$ws.Range("D17").Select()
